$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.006", "329.00") are preserved as text rather than being coerced
# into numbers, matching the inlineStr cell type used in the source file.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.666.92'
$ws.Range("E2").Value = '  +2.98%  '
$ws.Range("D3").Value = '2.004.23'
$ws.Range("E3").Value = '  +6.68%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("D5").Value = '329.00'
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("D7").Value = '0.4703'
$ws.Range("E7").Value = '  +2.36%  '
$ws.Range("E8").Value = '  +2.14%  '
$ws.Range("D9").Value = '46.98'
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("D10").Value = '0.07969'
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("E11").Value = '  +2.33%  '
$ws.Range("D12").Value = '22.85'
$ws.Range("E12").Value = '  +5.12%  '
$ws.Range("D13").Value = '2.003.48'
$ws.Range("E13").Value = '  +6.84%  '
$ws.Range("D14").Value = '7.279'
$ws.Range("E14").Value = '  +3.90%  '
$ws.Range("D15").Value = '5.903'
$ws.Range("E15").Value = '  +4.21%  '
$ws.Range("D16").Value = '0.07169'
$ws.Range("E16").Value = '  +3.38%  '
$ws.Range("D17").Value = '89.27'
$ws.Range("E17").Value = '  +1.07%  '
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").Value = '0.00001001'
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").Value = '17.44'
$ws.Range("E20").Value = '  +2.81%  '
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("D22").Value = '29.735.93'
$ws.Range("E22").Value = '  +3.19%  '
$ws.Range("D23").Value = '5.555'
$ws.Range("E23").Value = '  +5.45%  '
$ws.Range("D24").Value = '11.31'
$ws.Range("E24").Value = '  +3.30%  '
$ws.Range("D25").Value = '2.251.60'
$ws.Range("E25").Value = '  +7.19%  '
$ws.Range("D26").Value = '2.130'
$ws.Range("E26").Value = '  +2.08%  '
$ws.Range("D27").Value = '158.89'
$ws.Range("E27").Value = '  +2.23%  '
$ws.Range("D28").Value = '19.77'
$ws.Range("E28").Value = '  +2.57%  '
$ws.Range("D29").Value = '5.984'
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("D30").Value = '120.82'
$ws.Range("E30").Value = '  +2.89%  '
$ws.Range("D31").Value = '1.974'
$ws.Range("E31").Value = '  +2.27%  '
$ws.Range("D32").Value = '0.09490'
$ws.Range("E32").Value = '  +1.61%  '
$ws.Range("D33").Value = '0.9015'
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").Value = '5.309'
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("D35").Value = '1.347'
$ws.Range("E35").Value = '  +1.59%  '
$ws.Range("D36").Value = '3.196'
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("D37").Value = '0.05854'
$ws.Range("E37").Value = '  +1.67%  '
$ws.Range("D38").Value = '0.000003427'
$ws.Range("E38").Value = '  +111.75%  '
$ws.Range("D39").Value = '1.179'
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("D40").Value = '0.02141'
$ws.Range("E40").Value = '  +3.32%  '
$ws.Range("D41").Value = '7.930'
$ws.Range("E41").Value = '  +3.75%  '
$ws.Range("D42").Value = '0.5785'
$ws.Range("E42").Value = '  +2.34%  '
$ws.Range("D43").Value = '0.1827'
$ws.Range("E43").Value = '  +3.36%  '
$ws.Range("D44").Value = '9.896'
$ws.Range("E44").Value = '  +2.18%  '
$ws.Range("D45").Value = '12.21'
$ws.Range("E45").Value = '  +3.22%  '
$ws.Range("D46").Value = '0.5405'
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '2.185'
$ws.Range("E47").Value = '  -3.37%  '
$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").Value = '2.661'
$ws.Range("E48").Value = '  +5.31%  '
$ws.Range("D49").Value = '0.07001'
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").Value = '1.879'
$ws.Range("E50").Value = '  +1.65%  '
$ws.Range("D51").Value = '114.81'
$ws.Range("E51").Value = '  +1.76%  '

# Restore the cells to the workbook's default (unstyled) appearance now
# that the values have been written, so no stray number-format style is
# left attached to the cells.
$fmtRange.Style = "Normal"
